$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("modèle_de_feuille_de_données")

# The "Indemnités" header cell (L1) previously included a trailing
# ", 1/12 de la somme annuelle" clause that duplicated the note already
# shown under "Paiements spéciaux". Drop that clause while keeping the
# existing two-run rich text layout (bold title + regular note).
$rng = $ws.Range("L1")
$rng.Value = "Indemnités" + "`n" + "(Travail en équipes, le dimanche, de nuit et autres primes de pénibilité)"

$title = $rng.Characters(1, 10)
$title.Font.Name = "Arial"
$title.Font.Size = 9
$title.Font.Bold = $true

$note = $rng.Characters(11, 75)
$note.Font.Name = "Arial"
$note.Font.Size = 9
$note.Font.Bold = $false
$note.Font.Color = 0
